$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers
$ws.Range("F1").Value = "Baltimore"
$ws.Range("G1").Value = "Portland"

# Existing rows (2-6) get two new columns F and G
$ws.Range("F2").Value = 190
$ws.Range("G2").Value = 314

$ws.Range("F3").Value = 401
$ws.Range("G3").Value = 107

$ws.Range("F4").Value = 101
$ws.Range("G4").Value = 406

$ws.Range("F5").Value = 367
$ws.Range("G5").Value = 163

$ws.Range("F6").Value = 44
$ws.Range("G6").Value = 539

# New rows 7 (Baltimore) and 8 (Portland)
$ws.Range("A7").Value = 190
$ws.Range("B7").Value = 401
$ws.Range("C7").Value = 101
$ws.Range("D7").Value = 367
$ws.Range("E7").Value = 44
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 500

$ws.Range("A8").Value = 314
$ws.Range("B8").Value = 107
$ws.Range("C8").Value = 406
$ws.Range("D8").Value = 163
$ws.Range("E8").Value = 539
$ws.Range("F8").Value = 500
$ws.Range("G8").Value = 0

# Update selection to match target state
$ws.Range("J10").Select()
